$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("valueObject")

# A new "dynamicConditionList" field definition is being recorded at the
# bottom of the field-list table (rows 27-42). Insert a new row just above
# the trailing blank row (old row 43), copying row 42's formatting so the
# new row matches the rest of the table, then fill in the new field's data.
$ws.Rows("43:43").Insert()
$ws.Range("A42:F42").Copy()
$ws.Range("A43:F43").PasteSpecial(-4122)

# No. column keeps counting up from the row above.
$ws.Range("A43").Formula = "=A42+1"
$ws.Range("A42").Formula = "=A41+1"

# New field: dynamicConditionList
$ws.Range("B43").Value = "dynamicConditionList"
$ws.Range("C43").Value = "java.util.List<blanco.db.common.valueobject.BlancoDbDynamicConditionStructure>"
$ws.Range("D43").Value = "new java.util.ArrayList<blanco.db.common.valueobject.BlancoDbDynamicConditionStructure>()"
$ws.Range("E43").Value = "SQL動的条件定義のリスト"

# Re-assert inParameterList's type text on C32 (kept identical content).
$ws.Range("C32").Value = "java.util.List<blanco.dbmetadata.valueobject.BlancoDbMetaDataColumnStructure>"

# Update the view: select C32, matching where the author was working when
# the change was recorded.
$null = $ws.Range("C32").Select()
